$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B25").Value = 69

$ws.Range("I2").Select()
